$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "62.284.09"
$r.Style = "Normal"

$r = $ws.Range("E2")
$r.NumberFormat = "@"
$r.Value = "  -1.05%  "
$r.Style = "Normal"

$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "3.035.18"
$r.Style = "Normal"

$r = $ws.Range("E3")
$r.NumberFormat = "@"
$r.Value = "  -1.37%  "
$r.Style = "Normal"

$r = $ws.Range("E4")
$r.NumberFormat = "@"
$r.Value = "  +0.01%  "
$r.Style = "Normal"

$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "536.36"
$r.Style = "Normal"

$r = $ws.Range("E5")
$r.NumberFormat = "@"
$r.Value = "  -0.51%  "
$r.Style = "Normal"

$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "132.93"
$r.Style = "Normal"

$r = $ws.Range("E6")
$r.NumberFormat = "@"
$r.Value = "  +0.31%  "
$r.Style = "Normal"

$r = $ws.Range("E7")
$r.NumberFormat = "@"
$r.Value = "  -0.11%  "
$r.Style = "Normal"

$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "3.029.04"
$r.Style = "Normal"

$r = $ws.Range("E8")
$r.NumberFormat = "@"
$r.Value = "  -1.38%  "
$r.Style = "Normal"

$r = $ws.Range("E9")
$r.NumberFormat = "@"
$r.Value = "  +0.42%  "
$r.Style = "Normal"

$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "0.152"
$r.Style = "Normal"

$r = $ws.Range("E10")
$r.NumberFormat = "@"
$r.Value = "  -0.80%  "
$r.Style = "Normal"

$r = $ws.Range("E11")
$r.NumberFormat = "@"
$r.Value = "  +0.34%  "
$r.Style = "Normal"

$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "0.447"
$r.Style = "Normal"

$r = $ws.Range("E12")
$r.NumberFormat = "@"
$r.Value = "  -2.39%  "
$r.Style = "Normal"

$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "0.0000220"
$r.Style = "Normal"

$r = $ws.Range("E13")
$r.NumberFormat = "@"
$r.Value = "  -2.56%  "
$r.Style = "Normal"

$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "33.85"
$r.Style = "Normal"

$r = $ws.Range("E14")
$r.NumberFormat = "@"
$r.Value = "  -1.26%  "
$r.Style = "Normal"

$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "3.522.12"
$r.Style = "Normal"

$r = $ws.Range("E15")
$r.NumberFormat = "@"
$r.Value = "  -0.45%  "
$r.Style = "Normal"

$r = $ws.Range("E16")
$r.NumberFormat = "@"
$r.Value = "  +1.50%  "
$r.Style = "Normal"

$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "62.302.54"
$r.Style = "Normal"

$r = $ws.Range("E17")
$r.NumberFormat = "@"
$r.Value = "  -0.93%  "
$r.Style = "Normal"

$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "3.031.80"
$r.Style = "Normal"

$r = $ws.Range("E18")
$r.NumberFormat = "@"
$r.Value = "  -1.45%  "
$r.Style = "Normal"

$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "6.55"
$r.Style = "Normal"

$r = $ws.Range("E19")
$r.NumberFormat = "@"
$r.Value = "  -0.34%  "
$r.Style = "Normal"

$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "462.58"
$r.Style = "Normal"

$r = $ws.Range("E20")
$r.NumberFormat = "@"
$r.Value = "  -3.96%  "
$r.Style = "Normal"

$r = $ws.Range("E21")
$r.NumberFormat = "@"
$r.Value = "  +0.12%  "
$r.Style = "Normal"

$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "0.685"
$r.Style = "Normal"

$r = $ws.Range("E22")
$r.NumberFormat = "@"
$r.Value = "  -2.11%  "
$r.Style = "Normal"

$r = $ws.Range("E23")
$r.NumberFormat = "@"
$r.Value = "  -3.22%  "
$r.Style = "Normal"

$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "77.73"
$r.Style = "Normal"

$r = $ws.Range("E24")
$r.NumberFormat = "@"
$r.Value = "  -0.71%  "
$r.Style = "Normal"

$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "11.94"
$r.Style = "Normal"

$r = $ws.Range("E25")
$r.NumberFormat = "@"
$r.Value = "  -0.39%  "
$r.Style = "Normal"

$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "0.998"
$r.Style = "Normal"

$r = $ws.Range("E26")
$r.NumberFormat = "@"
$r.Value = "  -0.11%  "
$r.Style = "Normal"

$r = $ws.Range("E27")
$r.NumberFormat = "@"
$r.Value = "  -0.48%  "
$r.Style = "Normal"

$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "7.73"
$r.Style = "Normal"

$r = $ws.Range("E28")
$r.NumberFormat = "@"
$r.Value = "  -4.63%  "
$r.Style = "Normal"

$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "0.999"
$r.Style = "Normal"

$r = $ws.Range("E29")
$r.NumberFormat = "@"
$r.Value = "  -0.01%  "
$r.Style = "Normal"

$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "25.73"
$r.Style = "Normal"

$r = $ws.Range("E30")
$r.NumberFormat = "@"
$r.Value = "  -0.63%  "
$r.Style = "Normal"

$r = $ws.Range("B31")
$r.NumberFormat = "@"
$r.Value = "Mantle"
$r.Style = "Normal"

$r = $ws.Range("C31")
$r.NumberFormat = "@"
$r.Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$r.Style = "Normal"

$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "1.14"
$r.Style = "Normal"

$r = $ws.Range("E31")
$r.NumberFormat = "@"
$r.Value = "  +3.93%  "
$r.Style = "Normal"

$r = $ws.Range("B32")
$r.NumberFormat = "@"
$r.Value = "ImmutableX"
$r.Style = "Normal"

$r = $ws.Range("C32")
$r.NumberFormat = "@"
$r.Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$r.Style = "Normal"

$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "1.85"
$r.Style = "Normal"

$r = $ws.Range("E32")
$r.NumberFormat = "@"
$r.Value = "  -1.85%  "
$r.Style = "Normal"

$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "58.15"
$r.Style = "Normal"

$r = $ws.Range("E33")
$r.NumberFormat = "@"
$r.Value = "  -1.12%  "
$r.Style = "Normal"

$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "2.26"
$r.Style = "Normal"

$r = $ws.Range("E34")
$r.NumberFormat = "@"
$r.Value = "  -5.61%  "
$r.Style = "Normal"

$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "5.39"
$r.Style = "Normal"

$r = $ws.Range("E35")
$r.NumberFormat = "@"
$r.Value = "  +4.56%  "
$r.Style = "Normal"

$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "5.88"
$r.Style = "Normal"

$r = $ws.Range("E36")
$r.NumberFormat = "@"
$r.Value = "  -0.79%  "
$r.Style = "Normal"

$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "461.16"
$r.Style = "Normal"

$r = $ws.Range("E37")
$r.NumberFormat = "@"
$r.Value = "  -0.20%  "
$r.Style = "Normal"

$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "3.183.84"
$r.Style = "Normal"

$r = $ws.Range("E38")
$r.NumberFormat = "@"
$r.Value = "  +2.17%  "
$r.Style = "Normal"

$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "0.0389"
$r.Style = "Normal"

$r = $ws.Range("E39")
$r.NumberFormat = "@"
$r.Value = "  +0.22%  "
$r.Style = "Normal"

$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "0.0787"
$r.Style = "Normal"

$r = $ws.Range("E40")
$r.NumberFormat = "@"
$r.Value = "  +0.34%  "
$r.Style = "Normal"

$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "0.116"
$r.Style = "Normal"

$r = $ws.Range("E41")
$r.NumberFormat = "@"
$r.Value = "  +2.32%  "
$r.Style = "Normal"

$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "8.03"
$r.Style = "Normal"

$r = $ws.Range("E42")
$r.NumberFormat = "@"
$r.Value = "  +0.36%  "
$r.Style = "Normal"

$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "2.48"
$r.Style = "Normal"

$r = $ws.Range("E43")
$r.NumberFormat = "@"
$r.Value = "  -0.62%  "
$r.Style = "Normal"

$r = $ws.Range("E44")
$r.NumberFormat = "@"
$r.Value = "  +0.14%  "
$r.Style = "Normal"

$r = $ws.Range("E45")
$r.NumberFormat = "@"
$r.Value = "  -0.59%  "
$r.Style = "Normal"

$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "24.93"
$r.Style = "Normal"

$r = $ws.Range("E46")
$r.NumberFormat = "@"
$r.Value = "  +3.11%  "
$r.Style = "Normal"

$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "121.33"
$r.Style = "Normal"

$r = $ws.Range("E47")
$r.NumberFormat = "@"
$r.Value = "  +3.49%  "
$r.Style = "Normal"

$r = $ws.Range("E48")
$r.NumberFormat = "@"
$r.Value = "  +1.17%  "
$r.Style = "Normal"

$r = $ws.Range("E49")
$r.NumberFormat = "@"
$r.Value = "  -2.00%  "
$r.Style = "Normal"

$r = $ws.Range("E50")
$r.NumberFormat = "@"
$r.Value = "  +0.58%  "
$r.Style = "Normal"

$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "1.24"
$r.Style = "Normal"

$r = $ws.Range("E51")
$r.NumberFormat = "@"
$r.Value = "  +5.31%  "
$r.Style = "Normal"

Write-Host "Applied 91 cell updates to cryptos sheet"

